$d = $word.ActiveDocument

$d.Content.Find.Execute("79×31=", $true, $false, $false, $false, $false, $true, 1, $false, "61×77=", 2) | Out-Null
$d.Content.Find.Execute("69×51=", $true, $false, $false, $false, $false, $true, 1, $false, "71×74=", 2) | Out-Null
$d.Content.Find.Execute("96×57=", $true, $false, $false, $false, $false, $true, 1, $false, "71×28=", 2) | Out-Null
$d.Content.Find.Execute("79×96=", $true, $false, $false, $false, $false, $true, 1, $false, "61×39=", 2) | Out-Null
$d.Content.Find.Execute("95×23=", $true, $false, $false, $false, $false, $true, 1, $false, "47×84=", 2) | Out-Null
$d.Content.Find.Execute("37×16=", $true, $false, $false, $false, $false, $true, 1, $false, "92×66=", 2) | Out-Null
$d.Content.Find.Execute("48×21=", $true, $false, $false, $false, $false, $true, 1, $false, "46×64=", 2) | Out-Null
$d.Content.Find.Execute("29×63=", $true, $false, $false, $false, $false, $true, 1, $false, "38×24=", 2) | Out-Null
$d.Content.Find.Execute("65×54=", $true, $false, $false, $false, $false, $true, 1, $false, "17×60=", 2) | Out-Null
$d.Content.Find.Execute("21×95=", $true, $false, $false, $false, $false, $true, 1, $false, "47×46=", 2) | Out-Null
$d.Content.Find.Execute("55×97=", $true, $false, $false, $false, $false, $true, 1, $false, "32×56=", 2) | Out-Null
$d.Content.Find.Execute("66×92=", $true, $false, $false, $false, $false, $true, 1, $false, "75×43=", 2) | Out-Null
$d.Content.Find.Execute("17×99=", $true, $false, $false, $false, $false, $true, 1, $false, "93×88=", 2) | Out-Null
$d.Content.Find.Execute("46×48=", $true, $false, $false, $false, $false, $true, 1, $false, "26×30=", 2) | Out-Null
$d.Content.Find.Execute("52×92=", $true, $false, $false, $false, $false, $true, 1, $false, "54×17=", 2) | Out-Null
$d.Content.Find.Execute("27×87=", $true, $false, $false, $false, $false, $true, 1, $false, "66×52=", 2) | Out-Null
$d.Content.Find.Execute("11×64=", $true, $false, $false, $false, $false, $true, 1, $false, "31×11=", 2) | Out-Null
$d.Content.Find.Execute("69×70=", $true, $false, $false, $false, $false, $true, 1, $false, "38×95=", 2) | Out-Null
$d.Content.Find.Execute("75×75=", $true, $false, $false, $false, $false, $true, 1, $false, "17×74=", 2) | Out-Null
$d.Content.Find.Execute("41×98=", $true, $false, $false, $false, $false, $true, 1, $false, "40×44=", 2) | Out-Null
$d.Content.Find.Execute("64×28=", $true, $false, $false, $false, $false, $true, 1, $false, "78×45=", 2) | Out-Null
$d.Content.Find.Execute("77×68=", $true, $false, $false, $false, $false, $true, 1, $false, "73×31=", 2) | Out-Null
$d.Content.Find.Execute("85×25=", $true, $false, $false, $false, $false, $true, 1, $false, "71×59=", 2) | Out-Null
$d.Content.Find.Execute("21×82=", $true, $false, $false, $false, $false, $true, 1, $false, "41×83=", 2) | Out-Null
$d.Content.Find.Execute("61×36=", $true, $false, $false, $false, $false, $true, 1, $false, "61×32=", 2) | Out-Null
